# Updates cryptos list values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'65.500.30"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  +0.09%  "
$ws.Cells.Item(2, 5).Style = "Normal"

# Row 3
$ws.Cells.Item(3, 4).Value = "'3.579.54"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  +4.24%  "
$ws.Cells.Item(3, 5).Style = "Normal"

# Row 4
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  -0.13%  "
$ws.Cells.Item(4, 5).Style = "Normal"

# Row 5
$ws.Cells.Item(5, 4).Value = "'599.50"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  +1.13%  "
$ws.Cells.Item(5, 5).Style = "Normal"

# Row 6
$ws.Cells.Item(6, 4).Value = "'140.67"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  +4.02%  "
$ws.Cells.Item(6, 5).Style = "Normal"

# Row 7
$ws.Cells.Item(7, 4).Value = "'3.579.14"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  +4.24%  "
$ws.Cells.Item(7, 5).Style = "Normal"

# Row 8
$ws.Cells.Item(8, 5).Value = "'  +0.09%  "
$ws.Cells.Item(8, 5).Style = "Normal"

# Row 9
$ws.Cells.Item(9, 5).Value = "'  +1.57%  "
$ws.Cells.Item(9, 5).Style = "Normal"

# Row 10
$ws.Cells.Item(10, 5).Value = "'  +3.86%  "
$ws.Cells.Item(10, 5).Style = "Normal"

# Row 11
$ws.Cells.Item(11, 4).Value = "'7.18"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  -3.94%  "
$ws.Cells.Item(11, 5).Style = "Normal"

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.394"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  +4.45%  "
$ws.Cells.Item(12, 5).Style = "Normal"

# Row 13
$ws.Cells.Item(13, 4).Value = "'4.175.61"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  +4.03%  "
$ws.Cells.Item(13, 5).Style = "Normal"

# Row 14
$ws.Cells.Item(14, 4).Value = "'0.0000189"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  +4.90%  "
$ws.Cells.Item(14, 5).Style = "Normal"

# Row 15
$ws.Cells.Item(15, 4).Value = "'3.563.39"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  +2.55%  "
$ws.Cells.Item(15, 5).Style = "Normal"

# Row 16
$ws.Cells.Item(16, 4).Value = "'27.09"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +2.43%  "
$ws.Cells.Item(16, 5).Style = "Normal"

# Row 17
$ws.Cells.Item(17, 5).Value = "'  +1.55%  "
$ws.Cells.Item(17, 5).Style = "Normal"

# Row 18
$ws.Cells.Item(18, 4).Value = "'65.382.46"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  +0.01%  "
$ws.Cells.Item(18, 5).Style = "Normal"

# Row 19
$ws.Cells.Item(19, 4).Value = "'10.40"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  +4.73%  "
$ws.Cells.Item(19, 5).Style = "Normal"

# Row 20
$ws.Cells.Item(20, 4).Value = "'5.88"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  +2.10%  "
$ws.Cells.Item(20, 5).Style = "Normal"

# Row 21
$ws.Cells.Item(21, 4).Value = "'14.26"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  +3.87%  "
$ws.Cells.Item(21, 5).Style = "Normal"

# Row 22
$ws.Cells.Item(22, 4).Value = "'397.78"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  +1.71%  "
$ws.Cells.Item(22, 5).Style = "Normal"

# Row 23
$ws.Cells.Item(23, 4).Value = "'0.572"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  +4.88%  "
$ws.Cells.Item(23, 5).Style = "Normal"

# Row 24
$ws.Cells.Item(24, 2).Value = "'WrappedeETH"
$ws.Cells.Item(24, 2).Style = "Normal"
$ws.Cells.Item(24, 3).Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(24, 3).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'3.714.45"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  +3.88%  "
$ws.Cells.Item(24, 5).Style = "Normal"

# Row 25
$ws.Cells.Item(25, 2).Value = "'Litecoin"
$ws.Cells.Item(25, 2).Style = "Normal"
$ws.Cells.Item(25, 3).Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(25, 3).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'74.79"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  +2.50%  "
$ws.Cells.Item(25, 5).Style = "Normal"

# Row 26
$ws.Cells.Item(26, 5).Value = "'  +0.09%  "
$ws.Cells.Item(26, 5).Style = "Normal"

# Row 27
$ws.Cells.Item(27, 4).Value = "'0.0000117"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  +11.51%  "
$ws.Cells.Item(27, 5).Style = "Normal"

# Row 28
$ws.Cells.Item(28, 4).Value = "'7.87"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  +8.30%  "
$ws.Cells.Item(28, 5).Style = "Normal"

# Row 29
$ws.Cells.Item(29, 4).Value = "'0.997"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -0.37%  "
$ws.Cells.Item(29, 5).Style = "Normal"

# Row 30
$ws.Cells.Item(30, 4).Value = "'2.28"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  +1.17%  "
$ws.Cells.Item(30, 5).Style = "Normal"

# Row 31
$ws.Cells.Item(31, 4).Value = "'8.31"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  +2.18%  "
$ws.Cells.Item(31, 5).Style = "Normal"

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.590.45"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  +4.40%  "
$ws.Cells.Item(32, 5).Style = "Normal"

# Row 33
$ws.Cells.Item(33, 4).Value = "'24.15"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  +6.39%  "
$ws.Cells.Item(33, 5).Style = "Normal"

# Row 34
$ws.Cells.Item(34, 2).Value = "'Kaspa"
$ws.Cells.Item(34, 2).Style = "Normal"
$ws.Cells.Item(34, 3).Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(34, 3).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "'0.149"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  +2.44%  "
$ws.Cells.Item(34, 5).Style = "Normal"

# Row 35
$ws.Cells.Item(35, 2).Value = "'USDe"
$ws.Cells.Item(35, 2).Style = "Normal"
$ws.Cells.Item(35, 3).Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(35, 3).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "'1.00"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  +0.04%  "
$ws.Cells.Item(35, 5).Style = "Normal"

# Row 36
$ws.Cells.Item(36, 4).Value = "'1.28"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  +4.95%  "
$ws.Cells.Item(36, 5).Style = "Normal"

# Row 37
$ws.Cells.Item(37, 4).Value = "'7.08"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  +3.67%  "
$ws.Cells.Item(37, 5).Style = "Normal"

# Row 38
$ws.Cells.Item(38, 2).Value = "'ImmutableX"
$ws.Cells.Item(38, 2).Style = "Normal"
$ws.Cells.Item(38, 3).Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38, 3).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "'1.56"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  +2.40%  "
$ws.Cells.Item(38, 5).Style = "Normal"

# Row 39
$ws.Cells.Item(39, 2).Value = "'Monero"
$ws.Cells.Item(39, 2).Style = "Normal"
$ws.Cells.Item(39, 3).Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(39, 3).Style = "Normal"
$ws.Cells.Item(39, 4).Value = "'167.51"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  -3.23%  "
$ws.Cells.Item(39, 5).Style = "Normal"

# Row 40
$ws.Cells.Item(40, 4).Value = "'5.02"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  +4.99%  "
$ws.Cells.Item(40, 5).Style = "Normal"

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.0806"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  +4.07%  "
$ws.Cells.Item(41, 5).Style = "Normal"

# Row 42
$ws.Cells.Item(42, 4).Value = "'0.833"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  +2.36%  "
$ws.Cells.Item(42, 5).Style = "Normal"

# Row 43
$ws.Cells.Item(43, 4).Value = "'26.77"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  +18.76%  "
$ws.Cells.Item(43, 5).Style = "Normal"

# Row 44
$ws.Cells.Item(44, 4).Value = "'42.89"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -1.43%  "
$ws.Cells.Item(44, 5).Style = "Normal"

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.999"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -0.12%  "
$ws.Cells.Item(45, 5).Style = "Normal"

# Row 46
$ws.Cells.Item(46, 4).Value = "'4.46"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  +0.91%  "
$ws.Cells.Item(46, 5).Style = "Normal"

# Row 47
$ws.Cells.Item(47, 4).Value = "'1.70"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  +5.28%  "
$ws.Cells.Item(47, 5).Style = "Normal"

# Row 48
$ws.Cells.Item(48, 4).Value = "'1.20"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  +8.89%  "
$ws.Cells.Item(48, 5).Style = "Normal"

# Row 49
$ws.Cells.Item(49, 2).Value = "'Maker"
$ws.Cells.Item(49, 2).Style = "Normal"
$ws.Cells.Item(49, 3).Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'2.429.41"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  +10.51%  "
$ws.Cells.Item(49, 5).Style = "Normal"

# Row 50
$ws.Cells.Item(50, 2).Value = "'Cosmos"
$ws.Cells.Item(50, 2).Style = "Normal"
$ws.Cells.Item(50, 3).Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(50, 3).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'6.83"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  +4.80%  "
$ws.Cells.Item(50, 5).Style = "Normal"

# Row 51
$ws.Cells.Item(51, 4).Value = "'2.15"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  +2.00%  "
$ws.Cells.Item(51, 5).Style = "Normal"
